# concept closure, runs again.
# Adds a "Test" package and a "Specimen" concept (closing the concept
# set), then leaves the "concepts" sheet active/selected.

$wb = $excel.ActiveWorkbook

# --- packages: add a new "Test" package row ---
$packages = $wb.Worksheets.Item("packages")
$packages.Activate()
$packages.Range("A2").Value = "Test"
$packages.Range("B2").Value = "Test summary"
$packages.Range("D2").Select()

# --- concepts: add a new "Specimen" concept row under package "Test" ---
$concepts = $wb.Worksheets.Item("concepts")
$concepts.Activate()
$concepts.Range("A2").Value = "Test"
$concepts.Range("B2").Value = "Specimen"
$concepts.Range("C2").Value = "Specimen summary"
$concepts.Range("D2").Value = "Specimen description"
$concepts.Range("E2").Value = "default:C:Thing"
$concepts.Range("F2").Value = "default:C:Thing, default:c:Thing"
$concepts.Range("F2").Select()
